$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = '"Clinical and immunological features of severe and moderate coronavirus disease 2019"'
$ws.Range("E2").Value = '[Guang%Chen%xref no email%3,  Di%Wu%xref no email%3,  Wei%Guo%xref no email%3,  Yong%Cao%xref no email%3,  Da%Huang%xref no email%3,  Hongwu%Wang%xref no email%3,  Tao%Wang%xref no email%0,  Xiaoyun%Zhang%xref no email%3,  Huilong%Chen%xref no email%3,  Haijing%Yu%xref no email%3,  Xiaoping%Zhang%xref no email%3,  Minxia%Zhang%xref no email%3,  Shiji%Wu%xref no email%3,  Jianxin%Song%xref no email%3,  Tao%Chen%xref no email%3,  Meifang%Han%xref no email%3,  Shusheng%Li%xref no email%3,  Xiaoping%Luo%xref no email%3,  Jianping%Zhao%xref no email%3,  Qin%Ning%xref no email%3]'
$ws.Range("F2").Value = '10.1172/JCI137244'
$ws.Range("G2").Value = 'CROSSREF'
$ws.Range("H2").Value = '2023-05-26'
$ws.Range("E3").Value = '[Tao%Chen%NULL%0,          Di%Wu%NULL%0,          Huilong%Chen%NULL%0,          Weiming%Yan%NULL%0,          Danlei%Yang%NULL%0,          Guang%Chen%NULL%0,          Ke%Ma%NULL%0,          Dong%Xu%NULL%0,          Haijing%Yu%NULL%0,          Hongwu%Wang%NULL%0,          Tao%Wang%NULL%0,          Wei%Guo%NULL%0,          Jia%Chen%NULL%0,          Chen%Ding%NULL%0,          Xiaoping%Zhang%NULL%0,          Jiaquan%Huang%NULL%0,          Meifang%Han%NULL%0,          Shusheng%Li%NULL%0,          Xiaoping%Luo%NULL%0,          Jianping%Zhao%NULL%0,          Qin%Ning%NULL%0]'
$ws.Range("E4").Value = '[Ruchong%Chen%NULL%0,          Wenhua%Liang%NULL%0,          Mei%Jiang%NULL%0,          Weijie%Guan%NULL%0,          Chen%Zhan%NULL%2,          Tao%Wang%NULL%0,          Chunli%Tang%NULL%2,          Ling%Sang%NULL%0,          Jiaxing%Liu%NULL%2,          Zhengyi%Ni%NULL%0,          Yu%Hu%NULL%0,          Lei%Liu%NULL%0,          Hong%Shan%NULL%0,          Chunliang%Lei%NULL%0,          Yixiang%Peng%NULL%0,          Li%Wei%NULL%0,          Yong%Liu%NULL%0,          Yahua%Hu%NULL%0,          Peng%Peng%NULL%0,          Jianming%Wang%NULL%0,          Jiyang%Liu%NULL%0,          Zhong%Chen%NULL%0,          Gang%Li%NULL%0,          Zhijian%Zheng%NULL%0,          Shaoqin%Qiu%NULL%0,          Jie%Luo%NULL%0,          Changjiang%Ye%NULL%0,          Shaoyong%Zhu%NULL%0,          Xiaoqing%Liu%NULL%0,          Linling%Cheng%NULL%2,          Feng%Ye%NULL%0,          Jinping%Zheng%NULL%0,          Nuofu%Zhang%NULL%0,          Yimin%Li%NULL%0,          Jianxing%He%NULL%0,          Shiyue%Li%lishiyue@188.com%0,          Nanshan%Zhong%NULL%0,          NULL%NULL%NULL%0]'
$ws.Range("E5").Value = '[H.%Hou%NULL%1,          B.%Zhang%NULL%1,          H.%Huang%NULL%1,          Y.%Luo%NULL%1,          S.%Wu%NULL%1,          G.%Tang%NULL%1,          W.%Liu%NULL%1,          L.%Mao%NULL%2,          L.%Mao%NULL%0,          F.%Wang%fengwang@tjh.tjmu.edu.cn%1,          Z.%Sun%zysun@tjh.tjmu.edu.cn%1]'
$ws.Range("E6").Value = '[Fang%Liu%NULL%1,          Lin%Li%NULL%0,          MengDa%Xu%NULL%1,          Juan%Wu%NULL%1,          Ding%Luo%NULL%1,          YuSi%Zhu%NULL%1,          BiXi%Li%NULL%1,          XiaoYang%Song%NULL%1,          Xiang%Zhou%zhouxiang188483@126.com%1]'
$ws.Range("E7").Value = '[Kai-Cai%Liu%NULL%1,          Ping%Xu%xuping1027@163.com%1,          Wei-Fu%Lv%weifulv@ustc.edu.cn%1,          Xiao-Hui%Qiu%NULL%1,          Jin-Long%Yao%NULL%1,          Jin-Feng%Gu%NULL%1,          Wei%Wei%NULL%2]'
$ws.Range("E8").Value = '[Fengjun%Liu%NULL%1,          Qi%Zhang%NULL%0,          Chao%Huang%NULL%1,          Chunzi%Shi%NULL%1,          Lin%Wang%NULL%0,          Nannan%Shi%NULL%0,          Cong%Fang%NULL%1,          Fei%Shan%NULL%0,          Xue%Mei%NULL%1,          Jing%Shi%NULL%0,          Fengxiang%Song%NULL%1,          Zhongcheng%Yang%NULL%1,          Zezhen%Ding%NULL%1,          Xiaoming%Su%NULL%1,          Hongzhou%Lu%NULL%0,          Tongyu%Zhu%NULL%0,          Zhiyong%Zhang%NULL%0,          Lei%Shi%NULL%1,          Yuxin%Shi%NULL%0]'
$ws.Range("E9").Value = '[Lei%Pan%NULL%1,          Mi%Mu%NULL%0,          Pengcheng%Yang%NULL%0,          Yu%Sun%NULL%1,          Runsheng%Wang%NULL%0,          Junhong%Yan%NULL%1,          Pibao%Li%NULL%1,          Baoguang%Hu%NULL%1,          Jing%Wang%NULL%0,          Chao%Hu%NULL%0,          Yuan%Jin%NULL%1,          Xun%Niu%NULL%1,          Rongyu%Ping%NULL%0,          Yingzhen%Du%NULL%1,          Tianzhi%Li%NULL%0,          Guogang%Xu%NULL%0,          Qinyong%Hu%NULL%0,          Lei%Tu%NULL%0]'
$ws.Range("E10").Value = '[Haiying%Sun%NULL%0,          Ruoqi%Ning%NULL%0,          Yu%Tao%NULL%0,          Chong%Yu%NULL%0,          Xiaoyan%Deng%NULL%0,          Caili%Zhao%NULL%0,          Silu%Meng%NULL%0,          Fangxu%Tang%89650793@qq.com%0,          Dong%Xu%89650793@qq.com%0,          Dong%Xu%89650793@qq.com%0]'
$ws.Range("E11").Value = '[Sijia%Tian%NULL%0,          Nan%Hu%NULL%0,          Jing%Lou%NULL%0,          Kun%Chen%NULL%0,          Xuqin%Kang%NULL%0,          Zhenjun%Xiang%NULL%0,          Hui%Chen%NULL%0,          Dali%Wang%NULL%0,          Ning%Liu%NULL%0,          Dong%Liu%NULL%0,          Gang%Chen%NULL%0,          Yongliang%Zhang%NULL%0,          Dou%Li%NULL%0,          Jianren%Li%NULL%0,          Huixin%Lian%NULL%0,          Shengmei%Niu%NULL%0,          Luxi%Zhang%NULL%0,          Jinjun%Zhang%NULL%0]'
$ws.Range("E12").Value = '[Kelvin Kai-Wang%To%NULL%0,          Owen Tak-Yin%Tsang%NULL%1,          Wai-Shing%Leung%NULL%1,          Anthony Raymond%Tam%NULL%1,          Tak-Chiu%Wu%NULL%1,          David Christopher%Lung%NULL%1,          Cyril Chik-Yan%Yip%NULL%0,          Jian-Piao%Cai%NULL%0,          Jacky Man-Chun%Chan%NULL%1,          Thomas Shiu-Hong%Chik%NULL%1,          Daphne Pui-Ling%Lau%NULL%1,          Chris Yau-Chung%Choi%NULL%1,          Lin-Lei%Chen%NULL%1,          Wan-Mui%Chan%NULL%0,          Kwok-Hung%Chan%NULL%0,          Jonathan Daniel%Ip%NULL%0,          Anthony Chin-Ki%Ng%NULL%1,          Rosana Wing-Shan%Poon%NULL%0,          Cui-Ting%Luo%NULL%1,          Vincent Chi-Chung%Cheng%NULL%0,          Jasper Fuk-Woo%Chan%NULL%1,          Ivan Fan-Ngai%Hung%NULL%1,          Zhiwei%Chen%NULL%1,          Honglin%Chen%NULL%0,          Kwok-Yung%Yuen%NULL%0]'
$ws.Range("E13").Value = '[Suxin%Wan%NULL%0,          Qingjie%Yi%NULL%1,          Shibing%Fan%NULL%1,          Jinglong%Lv%13608388377@163.com%0,          Xianxiang%Zhang%NULL%2,          Xianxiang%Zhang%NULL%0,          Lian%Guo%NULL%1,          Chunhui%Lang%NULL%0,          Qing%Xiao%NULL%1,          Kaihu%Xiao%NULL%2,          Zhengjun%Yi%NULL%1,          Mao%Qiang%NULL%1,          Jianglin%Xiang%NULL%1,          Bangshuo%Zhang%NULL%1,          Yongping%Chen%NULL%1,          Cailiang%Gao%NULL%1]'
$ws.Range("E14").Value = '[Lang%Wang%NULL%0,          Wenbo%He%NULL%0,          Xiaomei%Yu%NULL%0,          Dalong%Hu%NULL%0,          Mingwei%Bao%NULL%0,          Huafen%Liu%NULL%0,          Jiali%Zhou%NULL%0,          Hong%Jiang%NULL%0]'
$ws.Range("E15").Value = '[Yu-Huan%Xu%NULL%0,          Jing-Hui%Dong%NULL%0,          Wei-Min%An%NULL%0,          Xiao-Yan%Lv%15001008285@139.com%0,          Xiao-Ping%Yin%yinxiaoping78@sina.com%0,          Jian-Zeng%Zhang%NULL%0,          Li%Dong%NULL%0,          Xi%Ma%NULL%0,          Hong-Jie%Zhang%NULL%0,          Bu-Lang%Gao%browngao@163.com%0]'
$ws.Range("E16").Value = '[Fanfan%Zeng%NULL%1,          Chan%Dai%NULL%2,          Chan%Dai%NULL%0,          Pengcheng%Cai%NULL%1,          Jinbiao%Wang%NULL%1,          Lei%Xu%NULL%0,          Jianyu%Li%NULL%1,          Guoyun%Hu%NULL%1,          Zheng%Wang%zhengwang@hust.edu.cn%1,          Fang%Zheng%zhengfangtj@hust.edu.cn%1,          Lin%Wang%lin_wang@hust.edu.cn%0]'
$ws.Range("E17").Value = '[Gemin%Zhang%NULL%0,          Jie%Zhang%945128911@qq.com%2,          Bowen%Wang%NULL%2,          Xionglin%Zhu%NULL%2,          Qiang%Wang%NULL%0,          Shiming%Qiu%NULL%2]'
$ws.Range("E18").Value = '[Gemin%Zhang%NULL%0,          Jie%Zhang%945128911@qq.com%0,          Bowen%Wang%NULL%0,          Xionglin%Zhu%NULL%0,          Qiang%Wang%NULL%0,          Shiming%Qiu%NULL%0]'
